$d = $word.ActiveDocument

# 1. Remove the "Meta description" paragraph that follows the title heading.
#    (It consists of an empty run, a bold "Meta description" run, and a plain
#    run with the rest of the sentence.)
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete() | Out-Null

# 2. Insert a new bold paragraph "Play Bounty Belles Free - Review and Pros &
#    Cons" right before the final paragraph of the document (the one that used
#    to hold the AI image-generation prompt).
$count = $d.Paragraphs.Count
$beforeLast = $d.Paragraphs.Item($count - 1)
$beforeLast.Range.InsertParagraphAfter() | Out-Null
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Bounty Belles Free - Review and Pros &amp; Cons</w:t></w:r></w:p>'
$newPara.Range.InsertXML($xml) | Out-Null

# 3. Replace the text of the last (italic) paragraph - the old AI image
#    prompt - with the meta description text, keeping its italic formatting.
$oldText = "Create an eye-catching feature image for the game Bounty Belles that fits with the Western theme and includes a happy Maya warrior wearing glasses. The image should be in a cartoon style and draw attention to the three determined girls who are the main characters of the game. The background should feature a desert landscape with a saloon, and the game logo and jackpot should be prominently displayed. Use colors that match the Western color scheme, such as brown, gold, and red. The image should convey the thrill of the Wild Bounty feature and the excitement of winning one of the three jackpots. Make sure to include the Maya warrior in a prominent position to give a unique twist to the Western theme."
$newText = "Check out our review of Bounty Belles and discover the pros and cons of this Western-themed slot game. Play for free now."
$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null
